# Scheduled market-data refresh: updates currentAveragePrice / leve price /
# profit columns (H:N) for a handful of leve rows across several sheets,
# reflecting newly pulled item prices. No formulas are involved - every
# touched cell stores a plain numeric literal.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 429.14285
$ws.Range("I92").Value = 429.14285
$ws.Range("K92").Value = 429.14285
$ws.Range("M92").Value = 818.85715

$ws.Range("H98").Value = 1149.1316
$ws.Range("I98").Value = 907.9091
$ws.Range("J98").Value = 2741.2
$ws.Range("K98").Value = 907.9091
$ws.Range("L98").Value = 2741.2
$ws.Range("M98").Value = 590.0909
$ws.Range("N98").Value = -5737.2

$ws.Range("H99").Value = 1094.2858
$ws.Range("I99").Value = 238.5
$ws.Range("K99").Value = 715.5
$ws.Range("M99").Value = 782.5

$ws.Range("H101").Value = 1784.7142
$ws.Range("J101").Value = 1778.3334
$ws.Range("L101").Value = 5335.0002
$ws.Range("N101").Value = -8579.0002

$ws.Range("H122").Value = 1149.1316
$ws.Range("I122").Value = 907.9091
$ws.Range("J122").Value = 2741.2
$ws.Range("K122").Value = 2723.7273
$ws.Range("L122").Value = 8223.599999999999
$ws.Range("M122").Value = -273.7273
$ws.Range("N122").Value = -13123.6

$ws.Range("H137").Value = 276585
$ws.Range("I137").Value = 513941.94
$ws.Range("J137").Value = 2711.6155
$ws.Range("K137").Value = 1541825.82
$ws.Range("L137").Value = 8134.8465
$ws.Range("M137").Value = -1539275.82
$ws.Range("N137").Value = -13234.8465

$ws.Range("H141").Value = 2761.7307
$ws.Range("I141").Value = 5202
$ws.Range("J141").Value = 2502.1277
$ws.Range("K141").Value = 15606
$ws.Range("L141").Value = 7506.3831
$ws.Range("M141").Value = -10426
$ws.Range("N141").Value = -17866.3831

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3976.76
$ws.Range("I32").Value = 3845.2122
$ws.Range("J32").Value = 17000
$ws.Range("K32").Value = 3845.2122
$ws.Range("L32").Value = 17000
$ws.Range("M32").Value = -3558.2122
$ws.Range("N32").Value = -17574

$ws.Range("H74").Value = 3175976.8
$ws.Range("I74").Value = 941.1836499999999
$ws.Range("K74").Value = 941.1836499999999
$ws.Range("M74").Value = -67.18364999999994

$ws.Range("H77").Value = 3175976.8
$ws.Range("I77").Value = 941.1836499999999
$ws.Range("K77").Value = 4705.91825
$ws.Range("M77").Value = -337.9182499999997

$ws.Range("H80").Value = 37655
$ws.Range("J80").Value = 37655
$ws.Range("L80").Value = 37655
$ws.Range("N80").Value = -39651

$ws.Range("H83").Value = 37655
$ws.Range("J83").Value = 37655
$ws.Range("L83").Value = 112965
$ws.Range("N83").Value = -122949

$ws.Range("H132").Value = 2019.9246
$ws.Range("I132").Value = 1498.2142
$ws.Range("J132").Value = 4011.9092
$ws.Range("K132").Value = 4494.642599999999
$ws.Range("L132").Value = 12035.7276
$ws.Range("M132").Value = -1964.642599999999
$ws.Range("N132").Value = -17095.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 27229.6
$ws.Range("J35").Value = 27229.6
$ws.Range("L35").Value = 27229.6
$ws.Range("N35").Value = -27849.6

$ws.Range("H82").Value = 26272.357
$ws.Range("I82").Value = 12959.25
$ws.Range("J82").Value = 31597.6
$ws.Range("K82").Value = 12959.25
$ws.Range("L82").Value = 31597.6
$ws.Range("M82").Value = -12576.25
$ws.Range("N82").Value = -32363.6

$ws.Range("H85").Value = 26272.357
$ws.Range("I85").Value = 12959.25
$ws.Range("J85").Value = 31597.6
$ws.Range("K85").Value = 12959.25
$ws.Range("L85").Value = 31597.6
$ws.Range("M85").Value = -11633.25
$ws.Range("N85").Value = -34249.6

$ws.Range("H102").Value = 11747.25
$ws.Range("I102").Value = 7329.6665
$ws.Range("J102").Value = 25000
$ws.Range("K102").Value = 7329.6665
$ws.Range("L102").Value = 25000
$ws.Range("M102").Value = -4084.6665
$ws.Range("N102").Value = -31490

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19670
$ws.Range("J41").Value = 23560
$ws.Range("L41").Value = 23560
$ws.Range("N41").Value = -24416

$ws.Range("H51").Value = 9103.799999999999
$ws.Range("J51").Value = 9103.799999999999
$ws.Range("L51").Value = 9103.799999999999
$ws.Range("N51").Value = -10575.8

$ws.Range("H61").Value = 9103.799999999999
$ws.Range("J61").Value = 9103.799999999999
$ws.Range("L61").Value = 9103.799999999999
$ws.Range("N61").Value = -9799.799999999999

$ws.Range("H109").Value = 11540
$ws.Range("J109").Value = 11540
$ws.Range("L109").Value = 11540
$ws.Range("N109").Value = -13620

$ws.Range("H134").Value = 6041.5806
$ws.Range("I134").Value = 8100.6816
$ws.Range("J134").Value = 1008.2222
$ws.Range("K134").Value = 24302.0448
$ws.Range("L134").Value = 3024.6666
$ws.Range("M134").Value = -21767.0448
$ws.Range("N134").Value = -8094.6666

$ws.Range("H135").Value = 49782.43
$ws.Range("J135").Value = 49782.43
$ws.Range("L135").Value = 49782.43
$ws.Range("N135").Value = -59922.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 30304984
$ws.Range("I118").Value = 55556388
$ws.Range("J118").Value = 3300.8
$ws.Range("K118").Value = 166669164
$ws.Range("L118").Value = 9902.400000000001
$ws.Range("M118").Value = -166667921
$ws.Range("N118").Value = -12388.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1679.68
$ws.Range("I113").Value = 1607.5883
$ws.Range("J113").Value = 1832.875
$ws.Range("K113").Value = 1607.5883
$ws.Range("L113").Value = 1832.875
$ws.Range("M113").Value = 562.4117000000001
$ws.Range("N113").Value = -6172.875

$ws.Range("H123").Value = 30026
$ws.Range("J123").Value = 30026
$ws.Range("L123").Value = 30026
$ws.Range("N123").Value = -34926

$ws.Range("H132").Value = 1916.275
$ws.Range("I132").Value = 1695.4062
$ws.Range("J132").Value = 2799.75
$ws.Range("K132").Value = 5086.2186
$ws.Range("L132").Value = 8399.25
$ws.Range("M132").Value = -2556.2186
$ws.Range("N132").Value = -13459.25

$ws.Range("H134").Value = 25499.834
$ws.Range("J134").Value = 25499.834
$ws.Range("L134").Value = 76499.50199999999
$ws.Range("N134").Value = -81569.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14290300
$ws.Range("I7").Value = 20002440
$ws.Range("J7").Value = 9952
$ws.Range("K7").Value = 20002440
$ws.Range("L7").Value = 9952
$ws.Range("M7").Value = -20002328
$ws.Range("N7").Value = -10176

$ws.Range("H22").Value = 1521.2778
$ws.Range("I22").Value = 1071.7142
$ws.Range("J22").Value = 1807.3636
$ws.Range("K22").Value = 1071.7142
$ws.Range("L22").Value = 1807.3636
$ws.Range("M22").Value = -776.7141999999999
$ws.Range("N22").Value = -2397.3636

$ws.Range("H27").Value = 1521.2778
$ws.Range("I27").Value = 1071.7142
$ws.Range("J27").Value = 1807.3636
$ws.Range("K27").Value = 1071.7142
$ws.Range("L27").Value = 1807.3636
$ws.Range("M27").Value = -964.7141999999999
$ws.Range("N27").Value = -2021.3636

$ws.Range("H41").Value = 9800
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H46").Value = 1767.3103
$ws.Range("I46").Value = 1650
$ws.Range("J46").Value = 1990.2
$ws.Range("K46").Value = 1650
$ws.Range("L46").Value = 1990.2
$ws.Range("M46").Value = -1462
$ws.Range("N46").Value = -2366.2

$ws.Range("H61").Value = 2649.5
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 3299
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 3299
$ws.Range("M61").Value = -1798
$ws.Range("N61").Value = -3703

$ws.Range("H113").Value = 2649.5
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 3299
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 3299
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -7639

$ws.Range("H126").Value = 14290300
$ws.Range("I126").Value = 20002440
$ws.Range("J126").Value = 9952
$ws.Range("K126").Value = 60007320
$ws.Range("L126").Value = 29856
$ws.Range("M126").Value = -60004850
$ws.Range("N126").Value = -34796

$ws.Range("H132").Value = 4990.478
$ws.Range("I132").Value = 5719.8667
$ws.Range("J132").Value = 3622.875
$ws.Range("K132").Value = 17159.6001
$ws.Range("L132").Value = 10868.625
$ws.Range("M132").Value = -14629.6001
$ws.Range("N132").Value = -15928.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1244.0204
$ws.Range("I132").Value = 1420.76
$ws.Range("J132").Value = 1059.9166
$ws.Range("K132").Value = 4262.28
$ws.Range("L132").Value = 3179.7498
$ws.Range("M132").Value = -1732.28
$ws.Range("N132").Value = -8239.7498
